# Update the "取得日時" (acquired timestamp) column on the "ランサーズ" sheet.
# All data rows (2-20) are stamped with the new run's timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-16 01:25:21"
$lastRow = 20

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
